# Slovenia Prva Liga - odds base update (03-04-2024)
# A new fixture row is inserted before the existing row 139 (id=137), pushing the
# three following rows down by one. The (new) row 140 also receives four updated
# odds values (N, P, R, S). Rows 141/142 keep their original data, just shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: create row 142 (brand-new row) -------------------------------
# Copy number formats only (reuse existing style indices, don't create new ones)
$ws.Range("A141").Copy()
$ws.Range("A142").PasteSpecial(-4122)
$ws.Range("E141").Copy()
$ws.Range("E142").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A142").Value = 140
$ws.Range("B142").Value = 6814434
$ws.Range("C142").Value = "Slovenia Prva Liga"
$ws.Range("D142").Value = "Slovenia Prva Liga"
$ws.Range("E142").Value = 45389.41666666666
$ws.Range("F142").Value = "NK Bravo"
$ws.Range("G142").Value = "NK Domzale"
$ws.Range("K142").Value = 1.833
$ws.Range("L142").Value = 3.25
$ws.Range("M142").Value = 4
$ws.Range("N142").Value = 1.833
$ws.Range("O142").Value = 3.25
$ws.Range("P142").Value = 4
$ws.Range("Q142").Value = -0.5
$ws.Range("R142").Value = 1.825
$ws.Range("S142").Value = 1.975
$ws.Range("T142").Value = 2.25
$ws.Range("U142").Value = 1.8
$ws.Range("V142").Value = 2
$ws.Range("W142").Value = 0
$ws.Range("X142").Value = 0
$ws.Range("Y142").Value = 0
$ws.Range("Z142").Value = 0
$ws.Range("AA142").Value = 0

# --- Step 2: row 141 <- old row 140 data (unchanged values) ---------------
$ws.Range("A141").Value = 139
$ws.Range("B141").Value = 6837117
$ws.Range("E141").Value = 45388.63541666666
$ws.Range("F141").Value = "NS Mura"
$ws.Range("G141").Value = "NK Celje"
$ws.Range("K141").Value = 5.25
$ws.Range("L141").Value = 4.2
$ws.Range("M141").Value = 1.5
$ws.Range("N141").Value = 5.25
$ws.Range("O141").Value = 4.2
$ws.Range("P141").Value = 1.5
$ws.Range("Q141").Value = 1
$ws.Range("R141").Value = 1.925
$ws.Range("S141").Value = 1.875
$ws.Range("T141").Value = 2.75
$ws.Range("U141").Value = 1.975
$ws.Range("V141").Value = 1.825

# --- Step 3: row 140 <- old row 139 data, with N/P/R/S odds updated -------
$ws.Range("A140").Value = 138
$ws.Range("B140").Value = 6814435
$ws.Range("E140").Value = 45388.52083333334
$ws.Range("F140").Value = "NK Radomlje"
$ws.Range("G140").Value = "FC Koper"
$ws.Range("K140").Value = 2.55
$ws.Range("L140").Value = 3.25
$ws.Range("M140").Value = 2.55
$ws.Range("N140").Value = 2.45
$ws.Range("O140").Value = 3.25
$ws.Range("P140").Value = 2.625
$ws.Range("Q140").Value = 0
$ws.Range("R140").Value = 1.825
$ws.Range("S140").Value = 1.975
$ws.Range("T140").Value = 2.25
$ws.Range("U140").Value = 1.8
$ws.Range("V140").Value = 2

# --- Step 4: row 139 <- brand-new fixture (duplicate of row 138's data) ---
$ws.Range("A139").Value = 137
$ws.Range("B139").Value = 8035687
$ws.Range("E139").Value = 45388.41666666666
$ws.Range("F139").Value = "NK Rogaska"
$ws.Range("G139").Value = "Olimpija Ljubljana"
$ws.Range("K139").Value = 6
$ws.Range("L139").Value = 4.333
$ws.Range("M139").Value = 1.45
$ws.Range("N139").Value = 5
$ws.Range("O139").Value = 4
$ws.Range("P139").Value = 1.571
$ws.Range("Q139").Value = 1
$ws.Range("R139").Value = 1.825
$ws.Range("S139").Value = 1.975
$ws.Range("T139").Value = 2.75
$ws.Range("U139").Value = 1.825
$ws.Range("V139").Value = 1.975
